$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (rows 4, 5, 9) ---
$ws.Rows.Item(4).RowHeight = 32.25
$ws.Rows.Item(5).RowHeight = 32.25
$ws.Rows.Item(9).RowHeight = 48

# --- Column width adjustments ---
# Column C gets an explicit width (closest achievable to 11.85546875 "character" width)
$ws.Columns.Item(3).ColumnWidth = 11
# Column F loses its "best fit" flag while keeping ~same width
$ws.Columns.Item(6).ColumnWidth = 12.6

# --- New content below the existing table (rows 64-84) ---

# Row 64
$ws.Range("C64").Value = "Example: Source Field = F5"

# Row 66 - header row (bold)
$ws.Range("D66").Value = "Lookup Table"
$ws.Range("D66").Font.Bold = $true
$ws.Range("E66").Value = "Lookup Field"
$ws.Range("E66").Font.Bold = $true
$ws.Range("F66").Value = "Lookup Condition"
$ws.Range("F66").Font.Bold = $true
$ws.Range("G66").Value = "Label"
$ws.Range("G66").Font.Bold = $true

# Row 67
$ws.Range("D67").Value = "Payee"
$ws.Range("E67").Value = "FirstName"
$ws.Range("F67").Value = "First_Split"
$ws.Range("G67").Value = "John"

# Row 68
$ws.Range("D68").Value = "Payee"
$ws.Range("E68").Value = "MiddleName"
$ws.Range("F68").Value = "Second_Split"
$ws.Range("G68").Value = "Vic"

# Row 69
$ws.Range("D69").Value = "Payee"
$ws.Range("E69").Value = "LastName"
$ws.Range("F69").Value = "Third_Split"

# Row 71 - rich text with a bold run in the middle
$full71 = "  1.1. Add a LookupFieldConditionEnum that would take care of this case scenario. It should initially contain these properties:"
$ws.Range("C71").Value = $full71
$ws.Range("C71").Characters(14, 24).Font.Bold = $true

# Row 73-76
$ws.Range("C73").Value = "First_Split"
$ws.Range("C74").Value = "Second_Split"
$ws.Range("C75").Value = "Third_Split"
$ws.Range("C76").Value = "Value"

# Row 78-81
$ws.Range("C78").Value = 'Note: The Split refers to the element on the array split from the original string value e.g.: "John Vic".'
$ws.Range("C79").Value = 'The string is delimeted by a space (" ") character.'
$ws.Range("C80").Value = 'The Value property refers to the actual field value "John Vic".'
$ws.Range("C81").Value = "A more appropriate example for the Value property usage is SourceField=F7."

# Row 83 - header row (bold)
$ws.Range("C83").Value = "Lookup Table"
$ws.Range("C83").Font.Bold = $true
$ws.Range("D83").Value = "Lookup Field"
$ws.Range("D83").Font.Bold = $true
$ws.Range("E83").Value = "Lookup Condition"
$ws.Range("E83").Font.Bold = $true
$ws.Range("F83").Value = "Label"
$ws.Range("F83").Font.Bold = $true

# Row 84
$ws.Range("C84").Value = "Department"
$ws.Range("D84").Value = "Name"
$ws.Range("E84").Value = "Value"
$ws.Range("F84").Value = "Marketing"

# --- View state: select all cells (closest achievable approximation) ---
$ws.Cells.Select()
